$wb = $excel.ActiveWorkbook

# Update Summary sheet
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = "Rahma Al Qassimi"
$wsSummary.Range("B4").Value = 2604.15
$wsSummary.Range("B6").Value = 397146
$wsSummary.Range("B7").Value = 278004
$wsSummary.Range("B8").Value = 119142
$wsSummary.Range("B9").Value = 1.43

# Update Assets sheet
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("B2").Value = "Luxury Car"
$wsAssets.Range("C2").Value = 395082
$wsAssets.Range("C3").Value = 2064
$wsAssets.Range("C4").Value = 397146

# Update Liabilities sheet
$wsLiabilities = $wb.Worksheets.Item("Liabilities")
$wsLiabilities.Range("C2").Value = 237049
$wsLiabilities.Range("D2").Value = 4939
$wsLiabilities.Range("E2").Value = 4
$wsLiabilities.Range("C3").Value = 40955
$wsLiabilities.Range("D3").Value = 2048
$wsLiabilities.Range("C4").Value = 278004

$wb.Save()
